$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string need to be forced
# to Text so COM does not silently coerce them to a Double (which would
# drop formatting such as trailing zeros, e.g. "1.00" -> 1).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '27.969.41'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '1.647.28'
$ws.Range('E3').Value = '  +1.75%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '213.40'
$ws.Range('E5').Value = '  +0.98%  '
Set-TextValue 'D6' '0.525'
$ws.Range('E6').Value = '  -0.15%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  -0.05%  '
Set-TextValue 'D8' '23.66'
$ws.Range('E8').Value = '  +3.64%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('E10').Value = '  +0.28%  '
Set-TextValue 'D11' '0.0872'
$ws.Range('D12').Value = '1.879.65'
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('D13').Value = '1.644.26'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('E15').Value = '  +2.47%  '
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').Value = '27.940.82'
$ws.Range('E17').Value = '  +1.50%  '
Set-TextValue 'D18' '232.12'
$ws.Range('E18').Value = '  +1.47%  '
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  +0.65%  '
Set-TextValue 'D21' '0.999'
$ws.Range('E21').Value = '  -0.10%  '
Set-TextValue 'D22' '10.70'
$ws.Range('E22').Value = '  +7.67%  '
Set-TextValue 'D23' '4.39'
$ws.Range('E23').Value = '  +2.36%  '
$ws.Range('E24').Value = '  +3.79%  '
Set-TextValue 'D25' '151.64'
$ws.Range('E25').Value = '  +1.80%  '
Set-TextValue 'D26' '6.93'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('E27').Value = '  +1.03%  '
$ws.Range('E28').Value = '  -0.13%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('E31').Value = '  +0.53%  '
Set-TextValue 'D32' '3.33'
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('D33').Value = '1.456.51'
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('E36').Value = '  -0.66%  '
Set-TextValue 'D37' '0.889'
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  +0.91%  '
Set-TextValue 'D40' '0.922'
$ws.Range('E40').Value = '  -2.42%  '
Set-TextValue 'D41' '69.53'
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +0.55%  '
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.789.53'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '1.78'
$ws.Range('E48').Value = '  +5.50%  '
Set-TextValue 'D49' '88.81'
$ws.Range('E49').Value = '  +2.97%  '
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0101'
$ws.Range('E51').Value = '  -4.15%  '
